# feat: light and dark dragon change
#
# Task_小游戏任务表 ("mini-game task table") edits:
#  1. The "获取土龙" (Get Earth Dragon) task's reward string drops the
#     "41" entry and gains a "40" entry in its place (E8).
#  2. Task #5, previously "跑酷" (Parkour), becomes "获取光暗龙"
#     (Get Light/Dark Dragon) with its own reward list + quest object guid.
#  3. The old "跑酷" task is preserved as the new task #6 in row 10.
#  4. The active selection/scroll position in the sheet view is updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix up 获取土龙's reward list (row 8) ---------------------------
$ws.Range("E8").Value = "39|1|1||40|1|1||41|1|1||42|1|1||43|1|1"

# --- 2. Row 9 turns into the new 获取光暗龙 task -------------------------
$ws.Range("B9").Value = "获取光暗龙"
$ws.Range("E9").Value = "44|1|1||45|1|1||46|1|1||47|1|1||48|1|1||49|1|1||50|1|1||51|1|1||52|1|1||53|1|1"
$ws.Range("F9").Value = "7949884C461020935235E5834F66108F"

# --- 3. Row 10 becomes the old 跑酷 task, shifted down as task #6 -------
$ws.Range("C10").Clear()
$ws.Range("E10").Clear()
$ws.Range("A10").Value = 6
$ws.Range("B10").Value = "跑酷"
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 0
$ws.Range("F10").Value = "0F362B364C669123BC0886AEC93884B0"

# --- 4. Update the saved selection / scroll position --------------------
$ws.Range("F13").Select() | Out-Null
